$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddEmployee")

# Update the employee data table (headers in row 1 stay the same)
$ws.Range("A2").Value = "John"
$ws.Range("B2").Value = "K"
$ws.Range("C2").Value = "Doe"

$ws.Range("A3").Value = "Katie"
$ws.Range("B3").Value = "K"
$ws.Range("C3").Value = "Ball"

$ws.Range("A4").Value = "Donald"
$ws.Range("B4").Value = "K"
$ws.Range("C4").Value = "Trump"

$ws.Range("A5").Value = "Mohammed"
$ws.Range("B5").Value = "K"
$ws.Range("C5").Value = "Salah"

# Update the active selection shown when the sheet is saved
$ws.Range("E11").Select()
